$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the batter name
$ws.Name = "Rahul Chahar"

# Set entire used range to Text format first so numeric-looking values
# (e.g. "6", "0.00", "-") are kept as text, matching the source data (t="str").
$ws.Range("A1:M6").NumberFormat = "@"

# Row 1
$ws.Cells.Item(1,1).Value = "matchNo"
$ws.Cells.Item(1,2).Value = "teamName"
$ws.Cells.Item(1,3).Value = "batterName"
$ws.Cells.Item(1,4).Value = "states"
$ws.Cells.Item(1,5).Value = "runs"
$ws.Cells.Item(1,6).Value = "balls"
$ws.Cells.Item(1,7).Value = "fours"
$ws.Cells.Item(1,8).Value = "sixes"
$ws.Cells.Item(1,9).Value = "sr"
$ws.Cells.Item(1,10).Value = "opponentTeamName"
$ws.Cells.Item(1,11).Value = "venue"
$ws.Cells.Item(1,12).Value = "date"
$ws.Cells.Item(1,13).Value = "result"

# Row 2
$ws.Cells.Item(2,1).Value = "13th"
$ws.Cells.Item(2,2).Value = "Mumbai Indians"
$ws.Cells.Item(2,3).Value = "Rahul Chahar"
$ws.Cells.Item(2,4).Value = "c †Pant b Avesh Khan"
$ws.Cells.Item(2,5).Value = "6"
$ws.Cells.Item(2,6).Value = "6"
$ws.Cells.Item(2,7).Value = "1"
$ws.Cells.Item(2,8).Value = "0"
$ws.Cells.Item(2,9).Value = "100.00"
$ws.Cells.Item(2,10).Value = "Delhi Capitals"
$ws.Cells.Item(2,11).Value = "Chennai"
$ws.Cells.Item(2,12).Value = "April 20"
$ws.Cells.Item(2,13).Value = "Capitals won by 6 wickets (with 5 balls remaining)"

# Row 3
$ws.Cells.Item(3,1).Value = "5th"
$ws.Cells.Item(3,2).Value = "Mumbai Indians"
$ws.Cells.Item(3,3).Value = "Rahul Chahar"
$ws.Cells.Item(3,4).Value = "c Shubman Gill b Russell"
$ws.Cells.Item(3,5).Value = "8"
$ws.Cells.Item(3,6).Value = "7"
$ws.Cells.Item(3,7).Value = "0"
$ws.Cells.Item(3,8).Value = "0"
$ws.Cells.Item(3,9).Value = "114.28"
$ws.Cells.Item(3,10).Value = "Kolkata Knight Riders"
$ws.Cells.Item(3,11).Value = "Chennai"
$ws.Cells.Item(3,12).Value = "April 13"
$ws.Cells.Item(3,13).Value = "Mumbai won by 10 runs"

# Row 4
$ws.Cells.Item(4,1).Value = "30th"
$ws.Cells.Item(4,2).Value = "Mumbai Indians"
$ws.Cells.Item(4,3).Value = "Rahul Chahar"
$ws.Cells.Item(4,4).Value = "c Raina b Bravo"
$ws.Cells.Item(4,5).Value = "0"
$ws.Cells.Item(4,6).Value = "1"
$ws.Cells.Item(4,7).Value = "0"
$ws.Cells.Item(4,8).Value = "0"
$ws.Cells.Item(4,9).Value = "0.00"
$ws.Cells.Item(4,10).Value = "Chennai Super Kings"
$ws.Cells.Item(4,11).Value = "Dubai (DSC)"
$ws.Cells.Item(4,12).Value = "September 19"
$ws.Cells.Item(4,13).Value = "Super Kings won by 20 runs"

# Row 5
$ws.Cells.Item(5,1).Value = "39th"
$ws.Cells.Item(5,2).Value = "Mumbai Indians"
$ws.Cells.Item(5,3).Value = "Rahul Chahar"
$ws.Cells.Item(5,4).Value = "lbw b Patel"
$ws.Cells.Item(5,5).Value = "0"
$ws.Cells.Item(5,6).Value = "1"
$ws.Cells.Item(5,7).Value = "0"
$ws.Cells.Item(5,8).Value = "0"
$ws.Cells.Item(5,9).Value = "0.00"
$ws.Cells.Item(5,10).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(5,11).Value = "Dubai (DSC)"
$ws.Cells.Item(5,12).Value = "September 26"
$ws.Cells.Item(5,13).Value = "RCB won by 54 runs"

# Row 6
$ws.Cells.Item(6,1).Value = "1st"
$ws.Cells.Item(6,2).Value = "Mumbai Indians"
$ws.Cells.Item(6,3).Value = "Rahul Chahar"
$ws.Cells.Item(6,4).Value = "run out (Kohli/†de Villiers)"
$ws.Cells.Item(6,5).Value = "0"
$ws.Cells.Item(6,6).Value = "0"
$ws.Cells.Item(6,7).Value = "0"
$ws.Cells.Item(6,8).Value = "0"
$ws.Cells.Item(6,9).Value = "-"
$ws.Cells.Item(6,10).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(6,11).Value = "Chennai"
$ws.Cells.Item(6,12).Value = "April 09"
$ws.Cells.Item(6,13).Value = "RCB won by 2 wickets"
